$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "Sean Steele - ssteele1812@gmail.com"
$ws.Range("B26").Select()
